$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Add new info for row 4 (Bdd sqlite 2)
$ws.Range("D4").Value = "lien bdd site ok"

# Add new info for row 5 (sur le site (avant cypress))
$ws.Range("D5").Value = "lien routes bdd ok mais formulaire bdd non"
$ws.Range("E5").Value = "liste des tickets à afficher"

# Update the existing cell C5 with the expanded text
$ws.Range("C5").Value = "pas de methode fetch => réécriture des fonctions. Deplus app,js inexistant"

# New rows for additional tests
$ws.Range("A6").Value = "tri asc"
$ws.Range("B6").Value = "ok"

$ws.Range("A7").Value = "tri desc"
$ws.Range("B7").Value = "ok"

$ws.Range("A8").Value = "autres tris"
$ws.Range("B8").Value = "[] (vide)"

$ws.Range("A9").Value = "del"
$ws.Range("B9").Value = "echec"
$ws.Range("C9").Value = "selection par id à revoir"

$ws.Range("A10").Value = "update"
$ws.Range("B10").Value = "à faire"

# Column widths: widen column C, add widths for new columns D and E
# (values chosen to land on the nearest representable width bucket for this engine)
$ws.Columns.Item(3).ColumnWidth = 67.16666666666667
$ws.Columns.Item(4).ColumnWidth = 39
$ws.Columns.Item(5).ColumnWidth = 22.833333333333332

# Update selection to match the final workbook state
$ws.Range("B14").Select()

Write-Host "edit applied"
